$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TestCaseID values in column A (rows 2-5): TC001-4 -> TC009-12
$ws.Range("A2").Value = "TC009_AddCustomer_Valid"
$ws.Range("A3").Value = "TC010_AddCustomer_CustomerId_Empty"
$ws.Range("A4").Value = "TC011_AddCustomer_CustomerName_Empty"
$ws.Range("A5").Value = "TC012_AddCustomer_Blank"

# Update the selected/active cell to reflect the new view state
$ws.Range("B11").Select()
